# Ads1299_defRegs.xlsx edit: "Stable 8kHz ExG. Will now work on acheiving 16k."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 changes ---
# F2: value 1 -> 0, and apply a new "Note" style (fill+border) with center alignment
$ws.Range("F2").Value = 0
$ws.Range("F2").Style = "Note"
$ws.Range("F2").HorizontalAlignment = -4108

# M2: new note referencing the new shared string
$ws.Range("M2").Value = "B5 is zero if you want to dictate SPI SCLK on nRF end"

# --- Row 14 changes ---
$ws.Range("J14").Value = 0

# --- Row 15 changes ---
$ws.Range("J15").Value = 0

# --- Selection moves from M14 to M3 ---
$ws.Range("M3").Select()
